$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at the top of the data block (row 32), pushing existing
# rows 32:102 down to 35:105. The new rows inherit formatting (incl. the
# date style in column D) from the row above, as Excel normally does.
$ws.Rows("32:34").Insert()

# Populate the 3 new rows with the new weekly price entries
# (Femacal de La Calera - Membrillo, date 2023-03-29 / serial 45014)

# Row 32: Especial
$ws.Range("A32").Value = 3
$ws.Range("B32").Value = "Femacal de La Calera"
$ws.Range("C32").Value = "Coquimbo"
$ws.Range("D32").Value = "2023-03-29"
$ws.Range("E32").Value = 5
$ws.Range("F32").Value = "Fruta"
$ws.Range("G32").Value = 100104
$ws.Range("H32").Value = "Frutos de pepita"
$ws.Range("I32").Value = 100104003
$ws.Range("J32").Value = "Membrillo"
$ws.Range("K32").Value = "Champion"
$ws.Range("L32").Value = "Especial"
$ws.Range("M32").Value = 75
$ws.Range("N32").Value = 16000
$ws.Range("O32").Value = 16000
$ws.Range("P32").Value = 16000
$ws.Range("Q32").Value = "`$/caja 18 kilos empedrada"
$ws.Range("R32").Value = "Región de O'Higgins"
$ws.Range("S32").Value = 889
$ws.Range("T32").Value = 18

# Row 33: Extra (doble especial)
$ws.Range("A33").Value = 3
$ws.Range("B33").Value = "Femacal de La Calera"
$ws.Range("C33").Value = "Coquimbo"
$ws.Range("D33").Value = "2023-03-29"
$ws.Range("E33").Value = 5
$ws.Range("F33").Value = "Fruta"
$ws.Range("G33").Value = 100104
$ws.Range("H33").Value = "Frutos de pepita"
$ws.Range("I33").Value = 100104003
$ws.Range("J33").Value = "Membrillo"
$ws.Range("K33").Value = "Champion"
$ws.Range("L33").Value = "Extra (doble especial)"
$ws.Range("M33").Value = 56
$ws.Range("N33").Value = 18000
$ws.Range("O33").Value = 18000
$ws.Range("P33").Value = 18000
$ws.Range("Q33").Value = "`$/caja 18 kilos empedrada"
$ws.Range("R33").Value = "Región de O'Higgins"
$ws.Range("S33").Value = 1000
$ws.Range("T33").Value = 18

# Row 34: Primera
$ws.Range("A34").Value = 3
$ws.Range("B34").Value = "Femacal de La Calera"
$ws.Range("C34").Value = "Coquimbo"
$ws.Range("D34").Value = "2023-03-29"
$ws.Range("E34").Value = 5
$ws.Range("F34").Value = "Fruta"
$ws.Range("G34").Value = 100104
$ws.Range("H34").Value = "Frutos de pepita"
$ws.Range("I34").Value = 100104003
$ws.Range("J34").Value = "Membrillo"
$ws.Range("K34").Value = "Champion"
$ws.Range("L34").Value = "Primera"
$ws.Range("M34").Value = 70
$ws.Range("N34").Value = 14000
$ws.Range("O34").Value = 14000
$ws.Range("P34").Value = 14000
$ws.Range("Q34").Value = "`$/caja 18 kilos empedrada"
$ws.Range("R34").Value = "Región de O'Higgins"
$ws.Range("S34").Value = 778
$ws.Range("T34").Value = 18
